# Generate Report for handback
# Update the Correspond Handoff Datetime (col D) and Correspond Handback
# DateTime (col G) for the "8320dd9d..." row (row 3) in both the zh-cn
# and de-de localization sheets, reflecting the latest handback pass.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 09:13:04"
$wsZhCn.Range("G3").Value = "2016-01-08 09:13:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 09:13:13"
$wsDeDe.Range("G3").Value = "2016-01-08 09:14:04"
